# "Add files via upload" -- the refreshed purchasing-plan workbook adds a
# new "unit_price" column (G) to Sheet1. Prices are filled in only for the
# finished-goods (FG) sku rows (rows 2-21); the raw-material (RM) rows
# (22-36) are left blank in the new column, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell.
$ws.Range("G1").Value = "unit_price"

# unit_price values, row-aligned with the existing FG sku rows (2-21).
$rows   = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)
$prices = @(900, 1100, 950, 1200, 1800, 3200, 1700, 3000, 2500, 4200, 1600, 1650, 1700, 600, 1100, 2200, 2800, 2600, 2400, 2800)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 7).Value = $prices[$i]
}

# Column A was widened (best-fit) to accommodate the longest sku_type/raw
# material label ("COFFEE_EXTRACT") once the sheet was revisited.
$ws.Columns.Item(1).AutoFit() | Out-Null

# The saved view scrolled back to the top of the sheet with the cursor
# left on I9 (just past the new data), replacing the old topLeftCell/E22
# selection state.
$ws.Range("A1").Select() | Out-Null
$ws.Range("I9").Select() | Out-Null
